$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028615175837273
$ws.Cells.Item(2, 4).Value = 1.033628575161154
$ws.Cells.Item(2, 5).Value = 1.028556961455731
$ws.Cells.Item(2, 6).Value = 1.039268064560952
$ws.Cells.Item(2, 9).Value = 1.035293062170012
$ws.Cells.Item(2, 10).Value = 1.03376651196347
$ws.Cells.Item(2, 11).Value = 1.036430376430139
$ws.Cells.Item(2, 12).Value = 1.031373432674532
$ws.Cells.Item(2, 13).Value = 1.042053735841488
$ws.Cells.Item(2, 14).Value = 1.035234578985145
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029456491743848
$ws.Cells.Item(3, 4).Value = 1.034278058167354
$ws.Cells.Item(3, 5).Value = 1.029268837868374
$ws.Cells.Item(3, 6).Value = 1.041433908877496
$ws.Cells.Item(3, 9).Value = 1.035537060312372
$ws.Cells.Item(3, 10).Value = 1.034249066124903
$ws.Cells.Item(3, 11).Value = 1.036889305101237
$ws.Cells.Item(3, 12).Value = 1.031893528756696
$ws.Cells.Item(3, 13).Value = 1.044026193057578
$ws.Cells.Item(3, 14).Value = 1.035717818428826
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030001156845887
$ws.Cells.Item(4, 4).Value = 1.034698427882088
$ws.Cells.Item(4, 5).Value = 1.029730086881902
$ws.Cells.Item(4, 6).Value = 1.042830280607586
$ws.Cells.Item(4, 9).Value = 1.035693588911551
$ws.Cells.Item(4, 10).Value = 1.034560939733587
$ws.Cells.Item(4, 11).Value = 1.037185688584822
$ws.Cells.Item(4, 12).Value = 1.032230011645004
$ws.Cells.Item(4, 13).Value = 1.045297094063035
$ws.Cells.Item(4, 14).Value = 1.036030134933807
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030230200166415
$ws.Cells.Item(5, 4).Value = 1.03487517755847
$ws.Cells.Item(5, 5).Value = 1.029924143080748
$ws.Cells.Item(5, 6).Value = 1.043416130498608
$ws.Cells.Item(5, 9).Value = 1.035759070280009
$ws.Cells.Item(5, 10).Value = 1.034691962491884
$ws.Cells.Item(5, 11).Value = 1.037310150926692
$ws.Cells.Item(5, 12).Value = 1.032371455827076
$ws.Cells.Item(5, 13).Value = 1.045830113443636
$ws.Cells.Item(5, 14).Value = 1.03616134375945
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03026866141523
$ws.Cells.Item(6, 4).Value = 1.034904856137715
$ws.Cells.Item(6, 5).Value = 1.029956734576919
$ws.Cells.Item(6, 6).Value = 1.043514428662064
$ws.Cells.Item(6, 9).Value = 1.03577004596328
$ws.Cells.Item(6, 10).Value = 1.034713956574487
$ws.Cells.Item(6, 11).Value = 1.037331040655822
$ws.Cells.Item(6, 12).Value = 1.032395204143608
$ws.Cells.Item(6, 13).Value = 1.045919536226666
$ws.Cells.Item(6, 14).Value = 1.036183369076172
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030004217074595
$ws.Cells.Item(7, 4).Value = 1.034700789517169
$ws.Cells.Item(7, 5).Value = 1.029732679293174
$ws.Cells.Item(7, 6).Value = 1.04283811338134
$ws.Cells.Item(7, 9).Value = 1.035694465145745
$ws.Cells.Item(7, 10).Value = 1.034562690815527
$ws.Cells.Item(7, 11).Value = 1.037187352195972
$ws.Cells.Item(7, 12).Value = 1.032231901681441
$ws.Cells.Item(7, 13).Value = 1.045304221240138
$ws.Cells.Item(7, 14).Value = 1.036031888502484
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028899444856651
$ws.Cells.Item(8, 4).Value = 1.033848047889583
$ws.Cells.Item(8, 5).Value = 1.028797415184082
$ws.Cells.Item(8, 6).Value = 1.040001090657068
$ws.Cells.Item(8, 9).Value = 1.035375803565898
$ws.Cells.Item(8, 10).Value = 1.033929670242501
$ws.Cells.Item(8, 11).Value = 1.036585592484727
$ws.Cells.Item(8, 12).Value = 1.031549212607747
$ws.Cells.Item(8, 13).Value = 1.042721474233479
$ws.Cells.Item(8, 14).Value = 1.035397968967636
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026954832428133
$ws.Cells.Item(9, 4).Value = 1.032346267080658
$ws.Cells.Item(9, 5).Value = 1.027154117801696
$ws.Cells.Item(9, 6).Value = 1.034961671165768
$ws.Cells.Item(9, 9).Value = 1.034803860327143
$ws.Cells.Item(9, 10).Value = 1.032811364658075
$ws.Cells.Item(9, 11).Value = 1.035520810053908
$ws.Cells.Item(9, 12).Value = 1.030345824604376
$ws.Cells.Item(9, 13).Value = 1.03812763983014
$ws.Cells.Item(9, 14).Value = 1.034278075261009
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025659874184878
$ws.Cells.Item(10, 4).Value = 1.031345667652077
$ws.Cells.Item(10, 5).Value = 1.026061821333089
$ws.Cells.Item(10, 6).Value = 1.031573098794509
$ws.Cells.Item(10, 9).Value = 1.034415490585623
$ws.Cells.Item(10, 10).Value = 1.032063913946222
$ws.Cells.Item(10, 11).Value = 1.034807978447521
$ws.Cells.Item(10, 12).Value = 1.029543307731717
$ws.Cells.Item(10, 13).Value = 1.03503458908863
$ws.Cells.Item(10, 14).Value = 1.033529563083411
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.025099488444958
$ws.Cells.Item(11, 4).Value = 1.030912537442871
$ws.Cells.Item(11, 5).Value = 1.025589619035054
$ws.Cells.Item(11, 6).Value = 1.030098510150858
$ws.Cells.Item(11, 9).Value = 1.034245628449161
$ws.Cells.Item(11, 10).Value = 1.031739804073306
$ws.Cells.Item(11, 11).Value = 1.034498604102578
$ws.Cells.Item(11, 12).Value = 1.029195749564083
$ws.Cells.Item(11, 13).Value = 1.033687630278928
$ws.Cells.Item(11, 14).Value = 1.0332049929373
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024891387274027
$ws.Cells.Item(12, 4).Value = 1.030751673986709
$ws.Cells.Item(12, 5).Value = 1.025414338122028
$ws.Cells.Item(12, 6).Value = 1.029549645993727
$ws.Cells.Item(12, 9).Value = 1.034182277921483
$ws.Cells.Item(12, 10).Value = 1.031619346075992
$ws.Cells.Item(12, 11).Value = 1.034383581016161
$ws.Cells.Item(12, 12).Value = 1.029066641562406
$ws.Cells.Item(12, 13).Value = 1.033186126067705
$ws.Cells.Item(12, 14).Value = 1.033084363875813
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024936023371956
$ws.Cells.Item(13, 4).Value = 1.030786178822412
$ws.Cells.Item(13, 5).Value = 1.02545193121432
$ws.Cells.Item(13, 6).Value = 1.029667431244056
$ws.Cells.Item(13, 9).Value = 1.034195878440319
$ws.Cells.Item(13, 10).Value = 1.031645187858937
$ws.Cells.Item(13, 11).Value = 1.034408258739349
$ws.Cells.Item(13, 12).Value = 1.02909433609287
$ws.Cells.Item(13, 13).Value = 1.033293754554714
$ws.Cells.Item(13, 14).Value = 1.033110242357054
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025082285696037
$ws.Cells.Item(14, 4).Value = 1.03089924000381
$ws.Cells.Item(14, 5).Value = 1.025575127890196
$ws.Cells.Item(14, 6).Value = 1.030053164266568
$ws.Cells.Item(14, 9).Value = 1.034240397108973
$ws.Cells.Item(14, 10).Value = 1.031729848389943
$ws.Cells.Item(14, 11).Value = 1.034489098457872
$ws.Cells.Item(14, 12).Value = 1.029185077648407
$ws.Cells.Item(14, 13).Value = 1.03364620013564
$ws.Cells.Item(14, 14).Value = 1.033195023115725
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025172409504688
$ws.Cells.Item(15, 4).Value = 1.030968903429206
$ws.Cells.Item(15, 5).Value = 1.025651048818922
$ws.Cells.Item(15, 6).Value = 1.030290675415966
$ws.Cells.Item(15, 9).Value = 1.034267792541702
$ws.Cells.Item(15, 10).Value = 1.031782001365538
$ws.Cells.Item(15, 11).Value = 1.03453889219312
$ws.Cells.Item(15, 12).Value = 1.029240985268095
$ws.Cells.Item(15, 13).Value = 1.033863195582698
$ws.Cells.Item(15, 14).Value = 1.033247250154527
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025697072299402
$ws.Cells.Item(16, 4).Value = 1.031374415939523
$ws.Cells.Item(16, 5).Value = 1.026093176092987
$ws.Cells.Item(16, 6).Value = 1.03167080510412
$ws.Cells.Item(16, 9).Value = 1.034426727932255
$ws.Cells.Item(16, 10).Value = 1.032085414365761
$ws.Cells.Item(16, 11).Value = 1.034828495530136
$ws.Cells.Item(16, 12).Value = 1.029566372698879
$ws.Cells.Item(16, 13).Value = 1.035123818229357
$ws.Cells.Item(16, 14).Value = 1.033551094036012
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.026026270373541
$ws.Cells.Item(17, 4).Value = 1.031628819605265
$ws.Cells.Item(17, 5).Value = 1.026370717148763
$ws.Cells.Item(17, 6).Value = 1.032534539525883
$ws.Cells.Item(17, 9).Value = 1.034525968970105
$ws.Cells.Item(17, 10).Value = 1.032275614339102
$ws.Cells.Item(17, 11).Value = 1.035009964686072
$ws.Cells.Item(17, 12).Value = 1.029770462796693
$ws.Cells.Item(17, 13).Value = 1.035912501693869
$ws.Cells.Item(17, 14).Value = 1.033741564115131
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02621831882149
$ws.Cells.Item(18, 4).Value = 1.031777222081457
$ws.Cells.Item(18, 5).Value = 1.026532676346314
$ws.Cells.Item(18, 6).Value = 1.033037638017348
$ws.Cells.Item(18, 9).Value = 1.034583691099199
$ws.Cells.Item(18, 10).Value = 1.03238651051769
$ws.Cells.Item(18, 11).Value = 1.035115743761308
$ws.Cells.Item(18, 12).Value = 1.029889499011537
$ws.Cells.Item(18, 13).Value = 1.036371791928313
$ws.Cells.Item(18, 14).Value = 1.033852617779013
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.02628380793847
$ws.Cells.Item(19, 4).Value = 1.031827825716688
$ws.Cells.Item(19, 5).Value = 1.026587912817053
$ws.Cells.Item(19, 6).Value = 1.033209063451232
$ws.Cells.Item(19, 9).Value = 1.034603345165856
$ws.Cells.Item(19, 10).Value = 1.0324243157499
$ws.Cells.Item(19, 11).Value = 1.03515180003413
$ws.Cells.Item(19, 12).Value = 1.029930086250029
$ws.Cells.Item(19, 13).Value = 1.03652827432225
$ws.Cells.Item(19, 14).Value = 1.03389047669899
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025990947145105
$ws.Cells.Item(20, 4).Value = 1.031601523142619
$ws.Cells.Item(20, 5).Value = 1.026340931942425
$ws.Cells.Item(20, 6).Value = 1.03244194199898
$ws.Cells.Item(20, 9).Value = 1.034515338266812
$ws.Cells.Item(20, 10).Value = 1.032255212260504
$ws.Cells.Item(20, 11).Value = 1.0349905018842
$ws.Cells.Item(20, 12).Value = 1.029748566489907
$ws.Cells.Item(20, 13).Value = 1.035827959747221
$ws.Cells.Item(20, 14).Value = 1.033721133063242
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025039213678114
$ws.Cells.Item(21, 4).Value = 1.030865945741488
$ws.Cells.Item(21, 5).Value = 1.025538846332688
$ws.Cells.Item(21, 6).Value = 1.029939607175599
$ws.Cells.Item(21, 9).Value = 1.034227294553862
$ws.Cells.Item(21, 10).Value = 1.031704919884806
$ws.Cells.Item(21, 11).Value = 1.034465296152105
$ws.Cells.Item(21, 12).Value = 1.029158356781803
$ws.Cells.Item(21, 13).Value = 1.03354244665576
$ws.Cells.Item(21, 14).Value = 1.033170059209251
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024441116381429
$ws.Cells.Item(22, 4).Value = 1.030403576954434
$ws.Cells.Item(22, 5).Value = 1.025035214890263
$ws.Cells.Item(22, 6).Value = 1.028359696832262
$ws.Cells.Item(22, 9).Value = 1.034044707196163
$ws.Cells.Item(22, 10).Value = 1.031358528957352
$ws.Cells.Item(22, 11).Value = 1.034134455202822
$ws.Cells.Item(22, 12).Value = 1.028787214001237
$ws.Cells.Item(22, 13).Value = 1.032098587471307
$ws.Cells.Item(22, 14).Value = 1.032823176366949
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024758151087704
$ws.Cells.Item(23, 4).Value = 1.030648676239592
$ws.Cells.Item(23, 5).Value = 1.025302135635278
$ws.Cells.Item(23, 6).Value = 1.029197874788301
$ws.Cells.Item(23, 9).Value = 1.034141641255516
$ws.Cells.Item(23, 10).Value = 1.031542195366922
$ws.Cells.Item(23, 11).Value = 1.034309899494122
$ws.Cells.Item(23, 12).Value = 1.028983968967174
$ws.Cells.Item(23, 13).Value = 1.032864667167981
$ws.Cells.Item(23, 14).Value = 1.033007103603887
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.026006908090744
$ws.Cells.Item(24, 4).Value = 1.031613857198945
$ws.Cells.Item(24, 5).Value = 1.026354390366239
$ws.Cells.Item(24, 6).Value = 1.03248378500699
$ws.Cells.Item(24, 9).Value = 1.034520142329298
$ws.Cells.Item(24, 10).Value = 1.032264431218624
$ws.Cells.Item(24, 11).Value = 1.034999296499481
$ws.Cells.Item(24, 12).Value = 1.029758460507539
$ws.Cells.Item(24, 13).Value = 1.035866162885177
$ws.Cells.Item(24, 14).Value = 1.033730365113339
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.027457306642392
$ws.Cells.Item(25, 4).Value = 1.032734410463461
$ws.Cells.Item(25, 5).Value = 1.027578381937228
$ws.Cells.Item(25, 6).Value = 1.036269443921689
$ws.Cells.Item(25, 9).Value = 1.034952963166455
$ws.Cells.Item(25, 10).Value = 1.033100810482681
$ws.Cells.Item(25, 11).Value = 1.035796605293148
$ws.Cells.Item(25, 12).Value = 1.030656975596545
$ws.Cells.Item(25, 13).Value = 1.039320498606976
$ws.Cells.Item(25, 14).Value = 1.034567932131885
